# Insert a new weekly price-report row for Choclo (Terminal La Palmera de
# La Serena) at sheet row 761, pushing the existing rows 761-810 down to
# 762-811 (dimension grows from A1:R810 to A1:R811).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 761..810 down by one, inheriting formatting from row 760
# (this is what gives the new D761 cell the date style s="2").
$ws.Range("A761").EntireRow.Insert()

# Populate the newly-opened row 761 with the new record's data.
$ws.Range("A761").Value = 8
$ws.Range("B761").Value = "Terminal La Palmera de La Serena"
$ws.Range("C761").Value = "Coquimbo"
$ws.Range("D761").Value = 45021
$ws.Range("E761").Value = 4
$ws.Range("F761").Value = 100112024
$ws.Range("G761").Value = "Choclo"
$ws.Range("H761").Value = "Dulce o Americano"
$ws.Range("I761").Value = "Primera"
$ws.Range("J761").Value = 17000
$ws.Range("K761").Value = 280
$ws.Range("L761").Value = 300
$ws.Range("M761").Value = 290
$ws.Range("N761").Value = "$/unidad"
$ws.Range("O761").Value = "Provincia de Limarí"
$ws.Range("P761").Value = 290
$ws.Range("Q761").Value = 1
$ws.Range("R761").Value = "Hortaliza"
